$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently sits right under the
#    document title (it is being relocated/re-purposed later in the document).
$metaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Meta description")) {
        $metaIdx = $i
        break
    }
}
if ($metaIdx -gt 0) {
    $d.Paragraphs($metaIdx).Range.Delete() | Out-Null
}

# 2. Insert a new bold paragraph with the page title text right before the
#    "Prompt: ..." paragraph near the end of the document.
$promptIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Prompt:")) {
        $promptIdx = $i
        break
    }
}
if ($promptIdx -gt 0) {
    $promptPara = $d.Paragraphs($promptIdx)
    $promptPara.Range.InsertParagraphBefore() | Out-Null

    $newPara = $d.Paragraphs($promptIdx)
    $newPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Buffalo Boost for Free - Slot Game Review</w:t></w:r></w:p>') | Out-Null

    # 3. Replace the body of the (now shifted) "Prompt: ..." paragraph with the
    #    meta-description text, keeping its existing italic run formatting.
    $oldText = "Prompt: Create a fun and eye-catching feature image for Buffalo Boost that showcases the game's setting and main character. The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be riding on the back of a buffalo through the magnificent and exotic Grand Canyon landscape. The warrior should be wearing sunglasses and have a wide smile on their face, exuding confidence and happiness. In the background, the beautiful landscape of the Grand Canyon should be visible, with foxes, snakes, and mountain goats jumping around in the distance. The image should be vibrant and colorful, with attention to detail that brings the image to life. The buffalo should be surrounded by golden light, indicating its importance in the game, and the Spinmatic logo should be visible in the corner of the image. This feature image will attract players to the game and give them a glimpse into the exciting and unique world of Buffalo Boost."
    $newText = "Read our review of Buffalo Boost slot game and play for free. Learn about Collect and Buy Feature, pay lines, symbols, and gameplay experience."
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}
